$d = $word.ActiveDocument

# Locate the start of question 7 ("Przykladem przynaleznosci gruntu jest"),
# robustly via Find, then delete everything from there through the end of
# the document body content.
$findRange = $d.Content
$found = $findRange.Find.Execute("Przykładem przynależności gruntu jest:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the start of question 7 in the document."
}

# Walk back from the found range to the beginning of its paragraph so we
# capture the "7. " label/run that precedes the found text too.
$para = $findRange.Paragraphs.First
$startPos = $para.Range.Start

$deleteRange = $d.Range($startPos, $d.Content.End)
$deleteRange.Delete()

# Insert the replacement questions 7-17 as raw OOXML right after the
# remaining content (end of question 6).
$insAt = $d.Content.End
$insertionPoint = $d.Range($insAt, $insAt)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>7. Miejscem zamieszkania osoby fizycznej jest:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>a) adres administracyjny, pod którym osoba jest zameldowana;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) gmina, w której osoba przebywa z zamiarem stałego pobytu;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) powiat, w którym osoba przebywa z zamiarem stałego pobytu;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>d) miejscowość, w której osoba przebywa z zamiarem stałego pobytu.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>8. Od jakiej chwili osoba fizyczna, zgodnie z Kodeksem cywilnym, posiada zdolność prawną:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>a) od chwili urodzenia;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) od chwili ukończenia 13 lat;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) od chwili ukończenia 18 lat;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>d) żadna z powyższych odpowiedzi nie jest prawidłowa.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>9. Kto zgodnie z Kodeksem cywilnym nie ma zdolności do czynności prawnych:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>a) tylko osoby, które są ubezwłasnowolnione całkowicie;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) tylko osoby, które nie ukończyły lat trzynastu;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) na przykład osoby, które nie ukończyły lat piętnastu;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>d) na przykład osoby ubezwłasnowolnione całkowicie.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>10. Które z poniższych jednostek należy zaliczyć do osób prawnych:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>a) Skarb Państwa;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) spółka cywilna;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) spółka komandytowo-cywilna;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>d) wszystkie powyższe odpowiedzi są nieprawidłowe.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>11. Osobami prawnymi nie są:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>a) Skarb Państwa;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) jednostki organizacyjne, którym przepisy szczególne przyznają osobowość prawną;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) spółka komandytowo-akcyjna i spółka z ograniczoną odpowiedzialnością;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>d) wszystkie powyższe odpowiedzi są nieprawidłowe.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>12. Typowa cecha subsydiarnej odpowiedzialności wspólnika spółki jawnej sprowadza się do tego że:</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>a) wszyscy wspólnicy odpowiadają w równych częściach za zobowiązania spółki;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) tylko niektórzy wspólnicy spółki odpowiadają w całości za jej zobowiązania;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>c) wierzyciel spółki może prowadzić egzekucję z majątku wspólników w przypadku, gdy egzekucja z majątku spółki okaże się bezskuteczna;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>d) wszystkie powyższe odpowiedzi są nieprawidłowe.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>13. Zgodnie z przepisami Konstytucji RP gmina to:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>a) ogólna jednostka samorządu terytorialnego;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) lokalna wspólnota samorządowa, którą z mocy prawa tworzą mieszkańcy określonego terytorium;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) regionalna wspólnota samorządowa, którą z mocy prawa tworzą mieszkańcy największej jednostki zasadniczego podziału terytorialnego kraju w celu wykonywania administracji publicznej;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>d) wszystkie powyższe odpowiedzi są nieprawidłowe.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">14. </w:t></w:r><w:r><w:t>Przedsiębiorcą jest:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>a) każda spółka prawa handlowego;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>b) osoba fizyczna prowadząca w swoim imieniu działalność zawodową;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) każda osoba prawna;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>d) wszystkie powyższe odpowiedzi są prawidłowe.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>15. Przedsiębiorca będący osoba fizyczną działa:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>a) pod firmą;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) pod nazwą;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) pod swoim imieniem i nazwiskiem, bez żadnych dodatków;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>d) wszystkie powyższe odpowiedzi są nieprawidłowe.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>16. Przedsiębiorca będący osoba fizyczną działa:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>a) pod firmą;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) pod nazwą;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>c) pod swoim imieniem i nazwiskiem, bez żadnych dodatków;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>d) wszystkie powyższe odpowiedzi są nieprawidłowe.</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>17. Firma osoby fizycznej:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>a) może zawierać jej pseudonim, ale tylko jeśli dodatkowo zawiera jej imię i nazwisko;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>b) nie może zawierać pseudonimu, gdyż byłoby to sprzeczne z zasadami współżycia społecznego;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>c) nie może być zbyta;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:tab/><w:t>d) wszystkie powyższe odpowiedzi są nieprawidłowe.</w:t></w:r></w:p><w:p/>
'@

$insertionPoint.InsertXML($newXml)

Write-Output "Replacement complete. Paragraphs now: $($d.Paragraphs.Count)"
